# MAS Report templates formatting update
#
# - Apply a custom currency number format ("$"#,##0.00) to the "Average
#   value" amount cells (B6:C6, merged) on the "5 Summary" sheet.
# - The amount cell no longer carries its placeholder 0 - clear it so the
#   cell is blank, ready for the real reported figure.
# - Leave the selection on C15, matching where the author last clicked
#   before saving the template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6:C6").NumberFormat = """$""#,##0.00"
$ws.Range("B6").Value = $null

[void]$ws.Range("C15").Select()
